# Add a new "user_mobile" column (R) to the add-to-cart order test data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("R1").Value = "user_mobile"

# Map of serial_key (column B) -> mobile number
$phoneByKey = @{
    "239060824HOZ" = "7896123569"
    "2390608248O4" = "7896324589"
}

for ($r = 2; $r -le 22; $r++) {
    $key = $ws.Cells.Item($r, 2).Value2
    $phone = $phoneByKey[$key]
    $ws.Cells.Item($r, 18).Value = [double]$phone
}

$ws.Range("D6").Select()
